# SF POC Files and Changes
# Adds a new "SFLoginTest" worksheet, wires it into the Regression sheet,
# flips the existing tests' Run flag to NO (only the new SF test runs),
# and appends a few extra LoginTest rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New SFLoginTest sheet, appended after the last tab.
# ---------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$sf = $wb.Worksheets.Add($null, $wb.Worksheets.Item($sheetCount))
$sf.Name = "SFLoginTest"

$sf.Range("A1").Value = "TestID"
$sf.Range("B1").Value = "Issue"
$sf.Range("C1").Value = "Description"
$sf.Range("D1").Value = "UserName"
$sf.Range("E1").Value = "Password"
$sf.Range("F1").Value = "ErrorMessage"

$sf.Range("A2").Value = "TC01_SF_ValidLogin"
$sf.Range("B2").Value = 1243
$sf.Range("D2").Value = "dineshk.krishnamurthy@nbcuni.com.uat"
$sf.Range("E2").Value = "*Gtrc`$1234"
$sf.Range("C2").Value = "Login to SF with valid credentials"

$sf.Hyperlinks.Add($sf.Range("D2"), "https://test.salesforce.com", "", "", "dineshk.krishnamurthy@nbcuni.com.uat") | Out-Null

# ---------------------------------------------------------------------
# 2. Regression sheet: turn every existing test OFF, add DownloadTest's
#    missing Mode, and append a new row for the SFLoginTest suite.
# ---------------------------------------------------------------------
$reg = $wb.Worksheets.Item("Regression")
$reg.Range("B2").Value = "NO"
$reg.Range("B3").Value = "NO"
$reg.Range("B4").Value = "NO"
$reg.Range("B5").Value = "NO"
$reg.Range("B7").Value = "NO"
$reg.Range("C7").Value = "serial"
$reg.Range("B8").Value = "NO"

$reg.Range("A9").Value = "SFLoginTest"
$reg.Range("B9").Value = "YES"
$reg.Range("C9").Value = "parallel"

# Data validation: drop the single-cell B2 rule and extend the B3:B66
# rule so it also covers B2.
$reg.Range("B2").Validation.Delete()
$reg.Range("B3:B66").Validation.Delete()
$reg.Range("B2:B66").Validation.Add(3, 1, 3, "YES,NO")
$reg.Range("B2:B66").Validation.ErrorTitle = "Invalid Run Option"
$reg.Range("B2:B66").Validation.ErrorMessage = "Run column should be either 'Y' OR 'N'"

# ---------------------------------------------------------------------
# 3. LoginTest sheet: append three more scenarios (copies of the
#    existing rows with fresh TestIDs).
# ---------------------------------------------------------------------
$login = $wb.Worksheets.Item("LoginTest")
$login.Range("A5").Value = "TC04_ValidLogin"
$login.Range("B5").Value = 1243
$login.Range("C5").Value = "Login with valid credentials"
$login.Range("D5").Value = "JamesB"
$login.Range("E5").Value = "Test321"

$login.Range("A6").Value = "TC05_InValidLogin"
$login.Range("B6").Value = 1244
$login.Range("C6").Value = "Login with invalid credentials"
$login.Range("D6").Value = "JamesB"
$login.Range("E6").Value = "test3210"
$login.Range("F6").Value = "Incorrect user name or password."

$login.Range("A7").Value = "TC06_LoginCreateAccount"
$login.Range("B7").Value = 1245
$login.Range("C7").Value = "Login with valid credentials from Create Account"
$login.Range("D7").Value = "JamesB"
$login.Range("E7").Value = "Test321"

# ---------------------------------------------------------------------
# 4. Restore "Regression" as the active/selected tab (adding the new
#    sheet made it active by default).
# ---------------------------------------------------------------------
$reg.Activate()
$reg.Range("D10").Select() | Out-Null
$login.Range("A3").Select() | Out-Null
$reg.Activate()

Write-Output "SF POC changes applied"
